# Update a few imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.912
$ws.Range("B9").Value = 6.484999999999999
$ws.Range("B18").Value = 5.972
$ws.Range("B20").Value = 6.37
